# Corrected ISO abbreviation from YK to YT
# Renames the "Geodata_YK" and "UUID_YK" sheets to "Geodata_YT" / "UUID_YT"
# and updates the active sheet/selection state left behind by the edit.

$wb = $excel.ActiveWorkbook

# Rename the two mis-labelled sheets (Yukon's ISO 3166-2 code is YT, not YK).
$wb.Worksheets.Item("Geodata_YK").Name = "Geodata_YT"
$wb.Worksheets.Item("UUID_YK").Name = "UUID_YT"

# SubjectNTopic_TBS_QC: selection moved from B18 to A18 (sheet stays inactive).
$wsSubj = $wb.Worksheets.Item("SubjectNTopic_TBS_QC")
$wsSubj.Range("A18").Select()

# Geodata_YT becomes the active/selected tab, with a new cell selection (B27).
# This also moves UUID_YT's sheetView off "tabSelected" since only one sheet
# can be the active tab in the saved workbook view.
$wsGeo = $wb.Worksheets.Item("Geodata_YT")
$wsGeo.Activate()
$wsGeo.Range("B27").Select()
